$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 1).Value = '1331947'
$ws.Cells.Item(2, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331947'
$ws.Cells.Item(2, 3).Value = '[Impact Brazil] - Robot Software Solution Development intern'
$ws.Cells.Item(2, 4).Value = 'São Paulo, SP, Brasil'
$ws.Cells.Item(2, 5).Value = 'No'
$ws.Cells.Item(2, 6).Value = '4 applicants'
$ws.Cells.Item(2, 7).Value = '3 - 6 Months'
$ws.Cells.Item(2, 8).Value = 'XD4 ROBOTICS LTDA'

# Row 3
$ws.Cells.Item(3, 1).Value = '1331945'
$ws.Cells.Item(3, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331945'
$ws.Cells.Item(3, 3).Value = 'Business Development Officer'
$ws.Cells.Item(3, 4).Value = 'Amman, Jordan'
$ws.Cells.Item(3, 5).Value = 'No'
$ws.Cells.Item(3, 6).Value = '0 applicants'
$ws.Cells.Item(3, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(3, 8).Value = 'International TEFL Training Institute'

# Row 4
$ws.Cells.Item(4, 1).Value = '1331939'
$ws.Cells.Item(4, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331939'
$ws.Cells.Item(4, 3).Value = 'Taste Hungary | Junior Technical Support Engineer - French Speaker [EU Only]'
$ws.Cells.Item(4, 4).Value = 'Budapeste, Hungria'
$ws.Cells.Item(4, 5).Value = 'Yes'
$ws.Cells.Item(4, 6).Value = '4 applicants'
$ws.Cells.Item(4, 7).Value = '6 - 18 Months'
$ws.Cells.Item(4, 8).Value = 'EATON'

# Row 5
$ws.Cells.Item(5, 1).Value = '1331926'
$ws.Cells.Item(5, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331926'
$ws.Cells.Item(5, 3).Value = 'Production Engineer'
$ws.Cells.Item(5, 4).Value = 'İzmir, Türkiye'
$ws.Cells.Item(5, 5).Value = 'No'
$ws.Cells.Item(5, 6).Value = '2 applicants'
$ws.Cells.Item(5, 7).Value = '3 - 6 Months'
$ws.Cells.Item(5, 8).Value = 'Norm Fasteners'

# Row 6
$ws.Cells.Item(6, 1).Value = '1331915'
$ws.Cells.Item(6, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331915'
$ws.Cells.Item(6, 3).Value = 'Logistics Responsible'
$ws.Cells.Item(6, 4).Value = 'İzmir, Türkiye'
$ws.Cells.Item(6, 5).Value = 'No'
$ws.Cells.Item(6, 6).Value = '2 applicants'
$ws.Cells.Item(6, 7).Value = '3 - 6 Months'
$ws.Cells.Item(6, 8).Value = 'Norm Fasteners'

# Row 7
$ws.Cells.Item(7, 1).Value = '1331908'
$ws.Cells.Item(7, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331908'
$ws.Cells.Item(7, 3).Value = 'Social Media Manager'
$ws.Cells.Item(7, 4).Value = 'Kalamata 241 00, Greece'
$ws.Cells.Item(7, 5).Value = 'No'
$ws.Cells.Item(7, 6).Value = '1 applicant'
$ws.Cells.Item(7, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(7, 8).Value = 'Train Your Pulse'

# Row 8
$ws.Cells.Item(8, 1).Value = '1331907'
$ws.Cells.Item(8, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331907'
$ws.Cells.Item(8, 3).Value = 'Business Development Representative'
$ws.Cells.Item(8, 4).Value = 'Kalamata 241 00, Greece'
$ws.Cells.Item(8, 5).Value = 'No'
$ws.Cells.Item(8, 6).Value = '0 applicants'
$ws.Cells.Item(8, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(8, 8).Value = 'Train Your Pulse'

# Row 9
$ws.Cells.Item(9, 1).Value = '1331906'
$ws.Cells.Item(9, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331906'
$ws.Cells.Item(9, 3).Value = 'Quality Engineer'
$ws.Cells.Item(9, 4).Value = 'İzmir, Türkiye'
$ws.Cells.Item(9, 5).Value = 'No'
$ws.Cells.Item(9, 6).Value = '1 applicant'
$ws.Cells.Item(9, 7).Value = '3 - 6 Months'
$ws.Cells.Item(9, 8).Value = 'Norm Fasteners'

# Row 10
$ws.Cells.Item(10, 1).Value = '1331837'
$ws.Cells.Item(10, 2).Value = 'https://aiesec.org/opportunity/global-talent/1331837'
$ws.Cells.Item(10, 3).Value = 'Community Manager'
$ws.Cells.Item(10, 4).Value = 'Le Bardo, Tunisie'
$ws.Cells.Item(10, 5).Value = 'No'
$ws.Cells.Item(10, 6).Value = '0 applicants'
$ws.Cells.Item(10, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(10, 8).Value = 'Pcmet horizon'

# Row 11
$ws.Cells.Item(11, 1).Value = '1328856'
$ws.Cells.Item(11, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328856'
$ws.Cells.Item(11, 3).Value = 'Accelerate Romania | Digital Marketing & Community Manager'
$ws.Cells.Item(11, 4).Value = 'Iași, Romania'
$ws.Cells.Item(11, 5).Value = 'No'
$ws.Cells.Item(11, 6).Value = '27 applicants'
$ws.Cells.Item(11, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(11, 8).Value = 'Nouveaux'

# Row 12
$ws.Cells.Item(12, 1).Value = '1328441'
$ws.Cells.Item(12, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328441'
$ws.Cells.Item(12, 3).Value = 'Accelerate Romania | Italian-Speaking Purchasing & Negotiation Specialist (EU Citizenship Required)'
$ws.Cells.Item(12, 4).Value = 'Iași, Romania'
$ws.Cells.Item(12, 5).Value = 'No'
$ws.Cells.Item(12, 6).Value = '5 applicants'
$ws.Cells.Item(12, 7).Value = '6 - 18 Months'
$ws.Cells.Item(12, 8).Value = 'Veo Wordwide Services - Iași'

# Row 13
$ws.Cells.Item(13, 1).Value = '1328227'
$ws.Cells.Item(13, 2).Value = 'https://aiesec.org/opportunity/global-talent/1328227'
$ws.Cells.Item(13, 3).Value = 'Policy & Advocacy Intern'
$ws.Cells.Item(13, 4).Value = 'Hyderabad, Telangana, India'
$ws.Cells.Item(13, 5).Value = 'No'
$ws.Cells.Item(13, 6).Value = '10 applicants'
$ws.Cells.Item(13, 7).Value = '9 - 12 Weeks'
$ws.Cells.Item(13, 8).Value = 'Arunodhaya Trust'

# Row 14
$ws.Cells.Item(14, 1).Value = '1327904'
$ws.Cells.Item(14, 2).Value = 'https://aiesec.org/opportunity/global-talent/1327904'
$ws.Cells.Item(14, 3).Value = 'UX Research Trainee'
$ws.Cells.Item(14, 4).Value = 'Bruxelles, Belgio'
$ws.Cells.Item(14, 5).Value = 'No'
$ws.Cells.Item(14, 6).Value = '174 applicants'
$ws.Cells.Item(14, 7).Value = '6 - 18 Months'
$ws.Cells.Item(14, 8).Value = 'UCB'

# Highlight E4 (PREMIUM = Yes) with yellow fill
$ws.Cells.Item(4, 5).Interior.Color = 65535

# Column width adjustments
$ws.Columns.Item(3).ColumnWidth = 101.1666667
$ws.Columns.Item(4).ColumnWidth = 29.1666667
$ws.Columns.Item(6).ColumnWidth = 16.1666667
$ws.Columns.Item(8).ColumnWidth = 39.1666667

